$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- 1. Unmerge everything that currently exists so cells can be freely rewritten ---
$ws.Range("A1:A2").UnMerge() | Out-Null
$ws.Range("B1:B2").UnMerge() | Out-Null
$ws.Range("C1:C2").UnMerge() | Out-Null
$ws.Range("D1:H1").UnMerge() | Out-Null
$ws.Range("I1:I2").UnMerge() | Out-Null
$ws.Range("J1:J2").UnMerge() | Out-Null
$ws.Range("K1:K2").UnMerge() | Out-Null
$ws.Range("L1:L2").UnMerge() | Out-Null
$ws.Range("M1:M2").UnMerge() | Out-Null

# --- 2. Clear the old header values that are moving / disappearing ---
$ws.Range("A1:M2").ClearContents() | Out-Null

# --- 3. Write the new header text layout ---
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Nama"
$ws.Range("C1").Value = "NIK"
$ws.Range("D1").Value = "Nilai"
$ws.Range("H1").Value = "Nama Mapel Peminatan"
$ws.Range("I1").Value = "Organisasi"
$ws.Range("J1").Value = "Jabatan"
$ws.Range("K1").Value = "Penghargaan"
$ws.Range("L1").Value = "Cita-cita"
$ws.Range("M1").Value = "Asal Sekolah"

$ws.Range("D2").Value = "Matematika"
$ws.Range("E2").Value = "Bahasa Indonesia"
$ws.Range("F2").Value = "Bahasa Inggris"
$ws.Range("G2").Value = "Mapel Peminatan"

# --- 4. Re-create the merges for the new layout ---
$ws.Range("A1:A2").Merge() | Out-Null
$ws.Range("B1:B2").Merge() | Out-Null
$ws.Range("C1:C2").Merge() | Out-Null
$ws.Range("D1:G1").Merge() | Out-Null
$ws.Range("H1:H2").Merge() | Out-Null
$ws.Range("I1:I2").Merge() | Out-Null
$ws.Range("J1:J2").Merge() | Out-Null
$ws.Range("K1:K2").Merge() | Out-Null
$ws.Range("L1:L2").Merge() | Out-Null
$ws.Range("M1:M2").Merge() | Out-Null

# --- 5. Alignment / formatting ---
# Center + vertical-center block: A1:C2, I1:M2
$ws.Range("A1:C2").HorizontalAlignment = $xlCenter
$ws.Range("A1:C2").VerticalAlignment = $xlCenter
$ws.Range("I1:M2").HorizontalAlignment = $xlCenter
$ws.Range("I1:M2").VerticalAlignment = $xlCenter

# "Nilai" merged header (D1:G1) - horizontal center only
$ws.Range("D1:G1").HorizontalAlignment = $xlCenter

# Sub headers under Nilai (D2:G2) - horizontal center only
$ws.Range("D2:G2").HorizontalAlignment = $xlCenter

# "Nama Mapel Peminatan" header (H1:H2) - center + vertical-center + wrap text
$ws.Range("H1:H2").HorizontalAlignment = $xlCenter
$ws.Range("H1:H2").VerticalAlignment = $xlCenter
$ws.Range("H1:H2").WrapText = $true

# --- 6. Row height ---
$ws.Rows.Item(1).RowHeight = 15

# --- 7. Column widths ---
# Columns D, E, F, M keep their original (untouched) widths - their content/position
# relative to the sheet's stored <col> width is unchanged, so they are left alone.
# Only the columns whose effective header width changed are touched here, using the
# closest values achievable through the ColumnWidth property.
$ws.Columns.Item(1).ColumnWidth = 6
$ws.Columns.Item(7).ColumnWidth = 16.17
$ws.Columns.Item(8).ColumnWidth = 14.17
$ws.Columns.Item(9).ColumnWidth = 15
$ws.Columns.Item(10).ColumnWidth = 10.67
$ws.Columns.Item(11).ColumnWidth = 11.67
$ws.Columns.Item(12).ColumnWidth = 10.17
$ws.Columns.Item(14).ColumnWidth = 16.67

# --- 8. Selection ---
$ws.Range("A1:A2").Select() | Out-Null

Write-Host "Edit complete"
